# Apply the "update database" edit to the Overview sheet:
#  - refresh the publish-date headers (G9/H9) for the new 1402-03-07 release
#  - refresh the latest (H) column figures with the newly published numbers
#  - H15 (impairment charge row) changes from a numeric -202 to the "-" placeholder

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 9: publish-date labels
$ws.Range("G9").Value = "1402-03-07 (8)"
$ws.Range("H9").Value = "1402-03-07 (2)"

# Latest period (column H) figures
$ws.Range("H12").Value = -17794
$ws.Range("H13").Value = 17218
$ws.Range("H14").Value = -1854
$ws.Range("H15").Value = "-"
$ws.Range("H16").Value = 213
$ws.Range("H17").Value = 15577
$ws.Range("H19").Value = 259
$ws.Range("H20").Value = 12970
$ws.Range("H21").Value = -2128
$ws.Range("H22").Value = 10842
$ws.Range("H24").Value = 10842
